$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EPBDS-13258: Property file is not read from dependent project.
# The nested-test result block (rows 36-41, "_res_.$Step5.$StepN") previously
# showed raw unresolved message keys (say.hello, say.hello.1, jar.say.hello,
# jar.say.hello.1) because the message bundle from the dependent project
# wasn't being read. After the fix these now resolve to the same messages
# as the corresponding top-level test rows (25-30).

$ws.Range("D36").Value = "Hello, from Project!"
$ws.Range("E36").Value = "Hello, from Project!"
$ws.Range("F36").Value = "Hello, from Project!"
$ws.Range("G36").Value = "Bonjour, for fr-CA!"

$ws.Range("D37").Value = "Hello, Parameter!"
$ws.Range("E37").Value = "Hello, Parameter!"
$ws.Range("F37").Value = "Hello, Parameter!"
$ws.Range("G37").Value = "Bonjour, Parameter!"
$ws.Range("D37:G37").Font.Name = "Arial"
$ws.Range("D37:G37").Font.Size = 12

$ws.Range("E38").Value = "Bonjour, from MessageBundle!"
$ws.Range("G38").Value = "Bonjour, from MessageBundle!"

$ws.Range("E39").Value = "Bonjour, Parameter!"
$ws.Range("G39").Value = "Bonjour, Parameter!"

$ws.Range("D40").Value = "Hello, from Jar!"
$ws.Range("E40").Value = "Hello, from Jar!"
$ws.Range("F40").Value = "Hello, from Jar!"
$ws.Range("G40").Value = "Hello, from Jar!"

$ws.Range("D41").Value = "Hello, Parameter!"
$ws.Range("E41").Value = "Hello, Parameter!"
$ws.Range("F41").Value = "Hello, Parameter!"
$ws.Range("G41").Value = "Hello, Parameter!"
